# SecretSantaNode/list.xlsx
# Changed random mapping logic: added santa emailid / child emailid columns
# to the "employees" rows, and random mapping (group/time/gift) columns,
# for every participant row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- First, stamp the same cell style already used by column A (s="1")
# onto every new cell we are about to populate, WITHOUT disturbing the
# existing style table (PasteSpecial formats-only reuses the style index).
$ws.Range("A1").Copy() | Out-Null
$ws.Range("B1:B20").PasteSpecial(-4122) | Out-Null
$ws.Range("D1:G20").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Data: for every row, the santa's name (B), the group (D), the
# shift/time (E), the city (F) and the company (G).
$names  = @{1="Mahikanth Nag"; 2="babanag"}
$groups = @{
    1="Night";  2="Day";   3="Day";   4="Day";   5="Day";   6="Night"
    7="Day";    8="Night"; 9="Night"; 10="Night";11="Day";  12="Night"
    13="Day";   14="Day";  15="Day";  16="Night";17="Night";18="Day"
    19="Night"; 20="Night"
}
$shifts = @{
    1="Idera";  2="Idera";  3="Idera";  4="FedEx";   5="Cambium"
    6="Idera";  7="Cambium";8="Cambium";9="FedEx";   10="Idera"
    11="FedEx"; 12="FedEx"; 13="Idera"; 14="Cambium";15="Cambium"
    16="FedEx"; 17="FedEx"; 18="Cambium";19="FedEx"; 20="FedEx"
}

# --- Column A: employee email addresses (unchanged text; A1 is already
# correct and A2/A3 are simply re-asserted here for clarity/completeness).
$emails = @{2="babanag95@gmail.com"; 3="secretsanta.accolite@gmail.com"}

for ($r = 1; $r -le 20; $r++) {

    if ($emails.ContainsKey($r)) {
        $ws.Cells.Item($r, 1).Value = $emails[$r]
    }

    if ($names.ContainsKey($r)) {
        $ws.Cells.Item($r, 2).Value = $names[$r]
    } else {
        $ws.Cells.Item($r, 2).Value = "secretsanta"
    }

    $ws.Cells.Item($r, 4).Value = $groups[$r]
    $ws.Cells.Item($r, 5).Value = $shifts[$r]
    $ws.Cells.Item($r, 6).Value = "Hyderabad"
    $ws.Cells.Item($r, 7).Value = "Accolite"
}

Write-Host "Updated employee / mapping columns for rows 1-20"
